$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Table body changes (DataProvider table, rows 10-14)
#    New logic: row 10 generates a fresh, randomised hostname (to avoid
#    collisions between test runs); rows 11-14 reuse it via {PREVIOUS}.
# ---------------------------------------------------------------------------

# Row 11: add / true (unchanged) / {PREVIOUS} / pass / SUCCESS
$ws.Range("D11").Value = "{PREVIOUS}"
$ws.Range("E11").Value = "pass"
$ws.Range("F11").Value = "EPP_UNEXPECTED_COMMAND_SUCCESS"

# Row 12: add / true / {PREVIOUS} / fail (unchanged) / FAILURE (unchanged)
$ws.Range("C12").Formula = '="true"'
$ws.Range("C12").Copy()
$ws.Range("C12").PasteSpecial(-4163)
$ws.Range("D12").Value = "{PREVIOUS}"
$ws.Range("I12").ClearContents()

# Row 13: rem / false (unchanged) / {PREVIOUS} / pass (unchanged) / SUCCESS (unchanged)
$ws.Range("D13").Value = "{PREVIOUS}"

# Row 14: rem / false (unchanged) / {PREVIOUS} / fail (unchanged) / FAILURE (unchanged)
$ws.Range("D14").Value = "{PREVIOUS}"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Notes box (B3) - add a new bullet describing {PREVIOUS}
# ---------------------------------------------------------------------------
$ws.Range("B3").Value = '* Values in {parentheses} are not literal values, they indicate that an appropriate value should be computed.
    * A value of `{EMPTY}` indicates an empty string.
    * A value of `{PREVIOUS}` indicates that the value computed from the same column in the preceding row should be used.
* Any row containing an "action" value of "add" and a passOrFail value of "pass" should result in an object property being added, which will be validated using an `<info>` command.
* The cell values for rows with action=rem do not occur in previous rows, so the objects should not have those properties, meaning a successful response is an error.'

# Row 10: add / false / <random hostname formula text> / fail / FAILURE
# (set last so its brand new shared string is appended after {PREVIOUS}/notes)
$ws.Range("C10").Formula = '="false"'
$ws.Range("C10").Copy()
$ws.Range("C10").PasteSpecial(-4163)
$ws.Range("D10").Value = '{"ns1.epp-16.rst." & RANDCHARS(18) & ".icann"}'
$ws.Range("E10").Value = "fail"
$ws.Range("F10").Value = "EPP_UNEXPECTED_COMMAND_FAILURE"
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Layout tweaks: wider hostname column, shorter notes row, new selection
# ---------------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 52.83
$ws.Rows(3).RowHeight = 96
$ws.Range("B3:F3").Select()
